$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update narrative text content (Administrar Personal -> Actor Admin wording) ---

# "Resumen" row: clarify actor wording ("El Aministrador del sistema" -> "El actor Admin")
$ws.Range("B14").Value = "El actor Admin puede modificar los datos, exceptuando los datos identificadores, de un Hostel Worker."

# "2a" alternate course: Hostel Worker doesn't exist
$ws.Range("C23").Value = "El Hostel Worker no existe en el sistema. Se avisa al Actor Admin que el Hostel Worker no existe."

# "5a" alternate course: data validity check
$ws.Range("C24").Value = "El sistema comprueba la validez de los datos, en caso de que los datos no sean correctos, se le avisa al Actor Admin."

# --- Update sheet view (scroll position / selection) ---
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("J25").Select()
